# edit.ps1
# Applies the "Update countries & provincias Spain" change to paises.xlsx
#
# Summary of changes:
#  1. The country ranking table (sorted descending by "Casos totales") is
#     refreshed with newer figures for China, Austria, Israel, Dinamarca,
#     Kuwait (no ranking change) plus Estonia, Letonia (whose updated totals
#     push them one place higher in the sort, displacing Argelia/Eslovenia/
#     Catar/Emiratos Arabes Unidos/Ucrania/Nueva Zelanda and Bulgaria by one
#     row respectively).
#  2. The "Datos actualizados..." timestamp footer is bumped from 09:55 to
#     10:20.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: China
$ws.Cells.Item(7, 1).Value = "China"
$ws.Cells.Item(7, 2).Value = 81589
$ws.Cells.Item(7, 3).Value = 35
$ws.Cells.Item(7, 4).Value = 76408
$ws.Cells.Item(7, 5).Value = 1863
$ws.Cells.Item(7, 6).Value = 429
$ws.Cells.Item(7, 7).Value = 6
$ws.Cells.Item(7, 8).Value = 3318

# Row 16: Austria
$ws.Cells.Item(16, 1).Value = "Austria"
$ws.Cells.Item(16, 2).Value = 10842
$ws.Cells.Item(16, 3).Value = 131
$ws.Cells.Item(16, 4).Value = 1436
$ws.Cells.Item(16, 5).Value = 9260
$ws.Cells.Item(16, 6).Value = 215
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 146

# Row 21: Israel
$ws.Cells.Item(21, 1).Value = "Israel"
$ws.Cells.Item(21, 2).Value = 6211
$ws.Cells.Item(21, 3).Value = 119
$ws.Cells.Item(21, 4).Value = 289
$ws.Cells.Item(21, 5).Value = 5891
$ws.Cells.Item(21, 6).Value = 107
$ws.Cells.Item(21, 7).Value = 5
$ws.Cells.Item(21, 8).Value = 31

# Row 28: Dinamarca
$ws.Cells.Item(28, 1).Value = "Dinamarca"
$ws.Cells.Item(28, 2).Value = 3355
$ws.Cells.Item(28, 3).Value = 248
$ws.Cells.Item(28, 4).Value = 894
$ws.Cells.Item(28, 5).Value = 2357
$ws.Cells.Item(28, 6).Value = 145
$ws.Cells.Item(28, 7).Value = 0
$ws.Cells.Item(28, 8).Value = 104

# Row 55: Estonia
$ws.Cells.Item(55, 1).Value = "Estonia"
$ws.Cells.Item(55, 2).Value = 858
$ws.Cells.Item(55, 3).Value = 79
$ws.Cells.Item(55, 4).Value = 45
$ws.Cells.Item(55, 5).Value = 802
$ws.Cells.Item(55, 6).Value = 16
$ws.Cells.Item(55, 7).Value = 6
$ws.Cells.Item(55, 8).Value = 11

# Row 56: Argelia
$ws.Cells.Item(56, 1).Value = "Argelia"
$ws.Cells.Item(56, 2).Value = 847
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 4).Value = 61
$ws.Cells.Item(56, 5).Value = 728
$ws.Cells.Item(56, 6).Value = 0
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 58

# Row 57: Eslovenia
$ws.Cells.Item(57, 1).Value = "Eslovenia"
$ws.Cells.Item(57, 2).Value = 841
$ws.Cells.Item(57, 3).Value = 0
$ws.Cells.Item(57, 4).Value = 10
$ws.Cells.Item(57, 5).Value = 816
$ws.Cells.Item(57, 6).Value = 31
$ws.Cells.Item(57, 7).Value = 0
$ws.Cells.Item(57, 8).Value = 15

# Row 58: Catar
$ws.Cells.Item(58, 1).Value = "Catar"
$ws.Cells.Item(58, 2).Value = 835
$ws.Cells.Item(58, 3).Value = 0
$ws.Cells.Item(58, 4).Value = 71
$ws.Cells.Item(58, 5).Value = 762
$ws.Cells.Item(58, 6).Value = 37
$ws.Cells.Item(58, 7).Value = 0
$ws.Cells.Item(58, 8).Value = 2

# Row 59: Emiratos Arabes Unidos
$ws.Cells.Item(59, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(59, 2).Value = 814
$ws.Cells.Item(59, 3).Value = 0
$ws.Cells.Item(59, 4).Value = 61
$ws.Cells.Item(59, 5).Value = 745
$ws.Cells.Item(59, 6).Value = 2
$ws.Cells.Item(59, 7).Value = 0
$ws.Cells.Item(59, 8).Value = 8

# Row 60: Ucrania
$ws.Cells.Item(60, 1).Value = "Ucrania"
$ws.Cells.Item(60, 2).Value = 804
$ws.Cells.Item(60, 3).Value = 10
$ws.Cells.Item(60, 4).Value = 13
$ws.Cells.Item(60, 5).Value = 771
$ws.Cells.Item(60, 6).Value = 0
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 20

# Row 61: Nueva Zelanda
$ws.Cells.Item(61, 1).Value = "Nueva Zelanda"
$ws.Cells.Item(61, 2).Value = 797
$ws.Cells.Item(61, 3).Value = 89
$ws.Cells.Item(61, 4).Value = 92
$ws.Cells.Item(61, 5).Value = 704
$ws.Cells.Item(61, 6).Value = 2
$ws.Cells.Item(61, 7).Value = 0
$ws.Cells.Item(61, 8).Value = 1

# Row 73: Letonia
$ws.Cells.Item(73, 1).Value = "Letonia"
$ws.Cells.Item(73, 2).Value = 458
$ws.Cells.Item(73, 3).Value = 12
$ws.Cells.Item(73, 4).Value = 1
$ws.Cells.Item(73, 5).Value = 457
$ws.Cells.Item(73, 6).Value = 3
$ws.Cells.Item(73, 7).Value = 0
$ws.Cells.Item(73, 8).Value = 0

# Row 74: Bulgaria
$ws.Cells.Item(74, 1).Value = "Bulgaria"
$ws.Cells.Item(74, 2).Value = 449
$ws.Cells.Item(74, 3).Value = 27
$ws.Cells.Item(74, 4).Value = 25
$ws.Cells.Item(74, 5).Value = 414
$ws.Cells.Item(74, 6).Value = 17
$ws.Cells.Item(74, 7).Value = 0
$ws.Cells.Item(74, 8).Value = 10

# Row 86: Kuwait
$ws.Cells.Item(86, 1).Value = "Kuwait"
$ws.Cells.Item(86, 2).Value = 317
$ws.Cells.Item(86, 3).Value = 0
$ws.Cells.Item(86, 4).Value = 81
$ws.Cells.Item(86, 5).Value = 236
$ws.Cells.Item(86, 6).Value = 14
$ws.Cells.Item(86, 7).Value = 0
$ws.Cells.Item(86, 8).Value = 0

# Update the "last refreshed" timestamp footer
$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 10:20"
